# Workbook was touched by a newer Excel/Aspose round-trip ("aspose and excel
# working"): a new data row was appended on Sheet1 and the selection moved.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 3: A3 = "ieri", B3:D3 = "aa" (two new shared strings).
$ws.Range("A3").Value = "ieri"
$ws.Range("B3:D3").Value = "aa"

# Selection left on I6 when the file was saved.
$ws.Range("I6").Select() | Out-Null
